$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03096638061106205
$ws.Range("C2").Value = 0.012773372232913971
$ws.Range("D2").Value = 0.01003422960639
$ws.Range("E2").Value = 0.005914391949772835
$ws.Range("F2").Value = 0.000007880567864049226
$ws.Range("J2").Value = 0.1264667510986328
$ws.Range("K2").Value = 1.4497603178024292
